$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '26.953.63'
$ws.Range("E2").Value = '  -1.03%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.818.49'
$ws.Range("E3").Value = '  +0.01%  '

$ws.Range("E4").Value = '  -0.08%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '310.20'
$ws.Range("E5").Value = '  -0.95%  '

$ws.Range("E6").Value = '  -0.05%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4643'
$ws.Range("E7").Value = '  -0.39%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3708'
$ws.Range("E8").Value = '  -1.50%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07362'
$ws.Range("E9").Value = '  -0.42%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8734'
$ws.Range("E10").Value = '  +0.34%  '

$ws.Range("E11").Value = '  -0.55%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.847.59'
$ws.Range("E12").Value = '  +1.53%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.357'
$ws.Range("E13").Value = '  -0.85%  '

$ws.Range("B14").Value = 'TRON'
$ws.Range("C14").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.07080'
$ws.Range("E14").Value = '  -0.01%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.523'
$ws.Range("E15").Value = '  -2.27%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.67'
$ws.Range("E16").Value = '  -0.57%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.000'
$ws.Range("E17").Value = '  -0.20%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008734'
$ws.Range("E18").Value = '  -0.27%  '

$ws.Range("E19").Value = '  -0.01%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.74'
$ws.Range("E20").Value = '  -1.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '26.974.18'
$ws.Range("E21").Value = '  -1.02%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.328'
$ws.Range("E22").Value = '  +0.38%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.59'
$ws.Range("E23").Value = '  -3.08%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.078.99'
$ws.Range("E24").Value = '  +1.45%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.907'
$ws.Range("E25").Value = '  -1.70%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '152.07'
$ws.Range("E26").Value = '  +0.39%  '

$ws.Range("E27").Value = '  -0.56%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.150'
$ws.Range("E28").Value = '  -3.81%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.308'
$ws.Range("E29").Value = '  +0.21%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '115.76'
$ws.Range("E30").Value = '  -1.18%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08919'
$ws.Range("E31").Value = '  -0.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.7599'
$ws.Range("E32").Value = '  -2.72%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.158'
$ws.Range("E33").Value = '  -1.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.470'
$ws.Range("E34").Value = '  -1.10%  '

$ws.Range("E35").Value = '  -0.33%  '

$ws.Range("E36").Value = '  -0.03%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.095'
$ws.Range("E37").Value = '  -0.28%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01957'
$ws.Range("E38").Value = '  -0.48%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05259'
$ws.Range("E39").Value = '  +0.30%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.942'
$ws.Range("E40").Value = '  +2.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '7.268'
$ws.Range("E41").Value = '  -0.01%  '

$ws.Range("B42").Value = 'TheSandbox'
$ws.Range("C42").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5349'
$ws.Range("E42").Value = '  +0.86%  '

$ws.Range("B43").Value = 'RenderToken'
$ws.Range("C43").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.387'
$ws.Range("E43").Value = '  +0.77%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.1666'
$ws.Range("E44").Value = '  -1.33%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '8.461'
$ws.Range("E45").Value = '  -1.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4959'
$ws.Range("E46").Value = '  -1.82%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '10.37'
$ws.Range("E47").Value = '  -1.10%  '

$ws.Range("E48").Value = '  +0.86%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.000'
$ws.Range("E49").Value = '  -0.02%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '103.47'
$ws.Range("E50").Value = '  -1.88%  '

$ws.Range("E51").Value = '  -0.75%  '
